$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "310.90"
Set-TextValue "E2" "8.14%"
Set-TextValue "G2" "8"

# Row 3
Set-TextValue "D3" "31.91"
Set-TextValue "E3" "7.75%"
Set-TextValue "G3" "8"

# Row 4
Set-TextValue "D4" "5.351"
Set-TextValue "E4" "5.11%"
Set-TextValue "G4" "8"

# Row 5
Set-TextValue "D5" "0.07631"
Set-TextValue "E5" "14.03%"
Set-TextValue "G5" "8"

# Row 6
Set-TextValue "D6" "7.843"
Set-TextValue "E6" "6.90%"
Set-TextValue "G6" "8"

# Row 7
Set-TextValue "D7" "3.718"
Set-TextValue "E7" "9.08%"
Set-TextValue "G7" "8"

# Row 8
Set-TextValue "D8" "1.577"
Set-TextValue "E8" "15.85%"
Set-TextValue "G8" "8"

# Row 9
Set-TextValue "D9" "0.9233"
Set-TextValue "E9" "0.46%"
Set-TextValue "G9" "8"

# Row 10
Set-TextValue "D10" "0.01699"
Set-TextValue "E10" "2,522.13%"
Set-TextValue "G10" "8"

# Row 11
Set-TextValue "D11" "0.1721"
Set-TextValue "E11" "8.40%"
Set-TextValue "G11" "8"

# Row 12
Set-TextValue "D12" "0.07645"
Set-TextValue "E12" "13.49%"
Set-TextValue "G12" "8"

# Row 13
Set-TextValue "D13" "0.08153"
Set-TextValue "E13" "5.83%"
Set-TextValue "G13" "8"

# Row 14
Set-TextValue "D14" "0.03016"
Set-TextValue "E14" "2.77%"
Set-TextValue "G14" "8"

# Row 15
Set-TextValue "D15" "0.09882"
Set-TextValue "E15" "10.05%"
Set-TextValue "G15" "8"

# Row 16
Set-TextValue "D16" "0.001518"
Set-TextValue "E16" "-3.77%"
Set-TextValue "G16" "8"

# Row 17
Set-TextValue "D17" "0.04560"
Set-TextValue "G17" "8"

# Row 18
Set-TextValue "D18" "0.006220"
Set-TextValue "E18" "-0.54%"
Set-TextValue "G18" "8"

# Row 19
Set-TextValue "D19" "3.491"
Set-TextValue "E19" "1.52%"
Set-TextValue "G19" "8"

# Row 20
Set-TextValue "E20" "0.71%"
Set-TextValue "G20" "8"

# Row 21
Set-TextValue "D21" "0.3323"
Set-TextValue "E21" "3.42%"
Set-TextValue "G21" "8"

# Row 22
Set-TextValue "D22" "0.1336"
Set-TextValue "E22" "2.08%"
Set-TextValue "G22" "8"

# Row 23
Set-TextValue "D23" "4.198"
Set-TextValue "E23" "2.54%"
Set-TextValue "G23" "8"

# Row 24
Set-TextValue "D24" "0.1629"
Set-TextValue "E24" "3.00%"
Set-TextValue "G24" "8"

# Row 25
Set-TextValue "D25" "0.001224"
Set-TextValue "E25" "3.09%"
Set-TextValue "G25" "8"

# Row 26
Set-TextValue "D26" "0.004493"
Set-TextValue "E26" "9.23%"
Set-TextValue "G26" "8"

# Row 27
Set-TextValue "D27" "0.0001300"
Set-TextValue "E27" "8.35%"
Set-TextValue "G27" "8"

# Row 28
Set-TextValue "E28" "7.66%"
Set-TextValue "G28" "8"

# Row 29
Set-TextValue "G29" "8"

# Row 30
Set-TextValue "G30" "8"

# Row 31
Set-TextValue "G31" "8"

# Row 32
Set-TextValue "G32" "8"

# Row 33
Set-TextValue "G33" "8"

# Row 34
Set-TextValue "G34" "8"

# Row 35
Set-TextValue "G35" "8"

# Row 36
Set-TextValue "G36" "8"

# Row 37
Set-TextValue "G37" "8"

# Row 38
Set-TextValue "G38" "8"

# Row 39
Set-TextValue "G39" "8"

# Row 40
Set-TextValue "D40" "0.04599"
Set-TextValue "G40" "8"

# Row 41
Set-TextValue "D41" "0.007274"
Set-TextValue "E41" "8.47%"
Set-TextValue "G41" "8"

# Row 42
Set-TextValue "D42" "0.1367"
Set-TextValue "E42" "10.44%"
Set-TextValue "G42" "8"

# Row 43
Set-TextValue "E43" "0.91%"
Set-TextValue "G43" "8"

# Row 44
Set-TextValue "D44" "0.01414"
Set-TextValue "E44" "6.40%"
Set-TextValue "G44" "8"

# Row 45
Set-TextValue "D45" "0.00006122"
Set-TextValue "E45" "7.39%"
Set-TextValue "G45" "8"

# Row 46
Set-TextValue "G46" "8"

# Row 47
Set-TextValue "D47" "0.01300"
Set-TextValue "E47" "-0.48%"
Set-TextValue "G47" "8"

# Row 48
Set-TextValue "G48" "8"

# Row 49
Set-TextValue "G49" "8"

# Row 50
Set-TextValue "G50" "8"

# Row 51
Set-TextValue "G51" "8"
